$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G232").Value = 48
$ws.Range("L232").Value = 1

$ws.Range("C233").Value = 309
$ws.Range("G233").Value = 50

$ws.Range("C234").Value = 398
$ws.Range("G234").Value = 62
$ws.Range("M234").Value = 1

$ws.Range("C235").Value = 285
$ws.Range("D235").Value = 21
$ws.Range("G235").Value = 83
$ws.Range("L235").Value = 0
$ws.Range("M235").Value = 0

$ws.Range("C236").Value = 199
$ws.Range("D236").Value = 8
$ws.Range("G236").Value = 88
$ws.Range("L236").Value = 0
$ws.Range("M236").Value = 1

$ws.Range("C237").Value = 523
$ws.Range("D237").Value = 9
$ws.Range("G237").Value = 86
$ws.Range("I237").Value = 7
$ws.Range("L237").Value = 3
$ws.Range("M237").ClearContents()

$ws.Range("C238").Value = 347
$ws.Range("D238").Value = 6
$ws.Range("F238").Value = 2
$ws.Range("G238").Value = 86
$ws.Range("I238").Value = 4
$ws.Range("L238").Value = 0
$ws.Range("M238").Value = 2

$ws.Range("C239").Value = 36
$ws.Range("D239").Value = 0
$ws.Range("E239").Value = 8
$ws.Range("F239").Value = 3
$ws.Range("G239").Value = 85
$ws.Range("I239").Value = 0
$ws.Range("L239").Value = 1
$ws.Range("M239").Value = 0

# Update view: frozen pane top-left cell and active selection
$ws.Application.ActiveWindow.ScrollRow = 231
$ws.Range("A2").Select()
